# feat : make TreatmentRecords_api and related forms
#
# Adds a new API-spec row (row 11) describing the "치료 데이터 추가"
# (treatment data creation) endpoint, and updates the sheet's view state
# to match where the author had scrolled/selected next.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# New row: Index=10, Method=POST, URL=/api5/create_treatment/, description=치료 데이터 추가
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = "POST"
$ws.Cells.Item(11, 3).Value = "/api5/create_treatment/"
$ws.Cells.Item(11, 4).Value = "치료 데이터 추가"

# Match the author's resulting view: scrolled so column C is left-most visible,
# with the next empty row's Method cell (B12) selected ready for further entry.
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("B12").Select()
